$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so numeric-looking values
# (e.g. "1.014", "29.569.96") are preserved verbatim instead of being
# reinterpreted as numbers/scientific notation by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.569.96'
$ws.Range("E2").Value = '  -2.70%  '

$ws.Range("D3").Value = '2.004.79'
$ws.Range("E3").Value = '  -4.26%  '

$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +1.27%  '

$ws.Range("D5").Value = '330.04'
$ws.Range("E5").Value = '  -3.80%  '

$ws.Range("D6").Value = '1.013'
$ws.Range("E6").Value = '  +1.25%  '

$ws.Range("D7").Value = '0.5003'
$ws.Range("E7").Value = '  -4.54%  '

$ws.Range("D8").Value = '0.4231'
$ws.Range("E8").Value = '  -4.34%  '

$ws.Range("D9").Value = '54.70'
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("D10").Value = '0.09029'
$ws.Range("E10").Value = '  -3.00%  '

$ws.Range("D12").Value = '23.31'
$ws.Range("E12").Value = '  -6.04%  '

$ws.Range("D13").Value = '2.055.83'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").Value = '8.069'
$ws.Range("E14").Value = '  -6.25%  '

$ws.Range("D15").Value = '6.468'
$ws.Range("E15").Value = '  -6.24%  '

$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +1.23%  '

$ws.Range("D17").Value = '94.47'
$ws.Range("E17").Value = '  -6.74%  '

$ws.Range("E18").Value = '  -3.79%  '

$ws.Range("D19").Value = '0.06693'
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").Value = '19.67'
$ws.Range("E20").Value = '  -6.93%  '

$ws.Range("D21").Value = '1.014'
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("D22").Value = '5.980'
$ws.Range("E22").Value = '  -5.45%  '

$ws.Range("D23").Value = '29.641.04'
$ws.Range("E23").Value = '  -2.48%  '

$ws.Range("D24").Value = '12.02'
$ws.Range("E24").Value = '  -4.24%  '

$ws.Range("D25").Value = '2.307'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").Value = '158.64'
$ws.Range("E26").Value = '  -2.58%  '

$ws.Range("D27").Value = '20.77'
$ws.Range("E27").Value = '  -4.81%  '

$ws.Range("D28").Value = '6.340'
$ws.Range("E28").Value = '  -7.02%  '

$ws.Range("D29").Value = '2.298'
$ws.Range("E29").Value = '  -8.36%  '

$ws.Range("D30").Value = '128.31'
$ws.Range("E30").Value = '  -3.71%  '

$ws.Range("D31").Value = '1.058'
$ws.Range("E31").Value = '  -7.24%  '

$ws.Range("D32").Value = '0.09950'
$ws.Range("E32").Value = '  -4.89%  '

$ws.Range("D33").Value = '1.562'
$ws.Range("E33").Value = '  -7.15%  '

$ws.Range("D34").Value = '5.840'
$ws.Range("E34").Value = '  -6.64%  '

$ws.Range("E35").Value = '  -1.65%  '

$ws.Range("D36").Value = '0.02462'
$ws.Range("E36").Value = '  -6.59%  '

$ws.Range("D37").Value = '9.296'
$ws.Range("E37").Value = '  -8.52%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.06414'
$ws.Range("E38").Value = '  -6.23%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.303'
$ws.Range("E39").Value = '  -2.96%  '

$ws.Range("D40").Value = '0.6569'
$ws.Range("E40").Value = '  -6.03%  '

$ws.Range("D41").Value = '11.67'
$ws.Range("E41").Value = '  -7.07%  '

$ws.Range("E42").Value = '  -7.48%  '

$ws.Range("E43").Value = '  +1.32%  '

$ws.Range("D44").Value = '0.6350'
$ws.Range("E44").Value = '  -6.86%  '

$ws.Range("D45").Value = '13.45'
$ws.Range("E45").Value = '  -6.30%  '

$ws.Range("D46").Value = '2.199'
$ws.Range("E46").Value = '  -6.40%  '

$ws.Range("D47").Value = '1.303'
$ws.Range("E47").Value = '  -4.91%  '

$ws.Range("D48").Value = '3.515'
$ws.Range("E48").Value = '  -3.22%  '

$ws.Range("D49").Value = '0.00000000339'
$ws.Range("E49").Value = '  -3.79%  '

$ws.Range("D50").Value = '0.06989'
$ws.Range("E50").Value = '  -3.46%  '

$ws.Range("D51").Value = '1.130'
$ws.Range("E51").Value = '  -6.82%  '
